# Update _Trans-files with the new regions + add transport possibility
# from MAR to DKW.
#
# The "Fueltrade" sheet (internal codeName Sheet3) holds one row per
# trade-link technology: Reg1 | Reg2 | Comm | Comm1 | Comm2 | Tech | TradeLink.
# Append a new link from region MAR to DKW carrying commodity H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

$newRow = 44

$ws.Cells.Item($newRow, 2).Value = "MAR"                 # Reg1
$ws.Cells.Item($newRow, 3).Value = "DKW"                 # Reg2
$ws.Cells.Item($newRow, 4).Value = "H2"                  # Comm
$ws.Cells.Item($newRow, 5).Value = "H2"                  # Comm1
$ws.Cells.Item($newRow, 6).Value = "H2"                  # Comm2
$ws.Cells.Item($newRow, 7).Value = "TB_H2_MAR_DKW_01"    # Tech
$ws.Cells.Item($newRow, 8).Value = "B"                   # TradeLink

# The BI sheet becomes the active/selected tab (it was Fueltrade before).
$wsBI = $wb.Worksheets.Item("BI")
$wsBI.Select()
